# Improve random number generator.
# Update the "Price" and "Count" sample values on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.022477936010099
$ws.Range("C2").Value = 8

$ws.Range("B3").Value = 601.352605317418
$ws.Range("C3").Value = 89

$ws.Range("B4").Value = 967.955701969543
$ws.Range("C4").Value = 18
